$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells with known styles (row 3 is untouched by this edit):
# A3 -> style 1 (bold label), B3 -> style 2 (wrapped body), C3 -> style 3 (red wrapped body)

# --- Row 13 ---
$ws.Rows.Item(13).Clear()
$ws.Range("B3:C3").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B13").Value = '5840793 - Sérgio Schneider'
$ws.Range("C13").Value = '5840793 - Sérgio Schneider'
$ws.Rows.Item(13).AutoFit()

# --- Row 14 ---
$ws.Rows.Item(14).Clear()
$ws.Range("A3:C3").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Apresentar as principais técnicas ópticas para a medição de grandezas como comprimento, deslocamento e forma, com ênfase nas técnicas interferométricas a laser.'
$ws.Range("C14").Value = 'Apresentar as principais técnicas ópticas para a medição de grandezas como comprimento, deslocamento e forma, com ênfase nas técnicas interferométricas a laser.'
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15 ---
$ws.Rows.Item(15).Clear()
$ws.Range("A3:C3").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Present the main optical techniques for measuring quantities such as length, displacement and shape, with emphasis on laser interferometric techniques.'
$ws.Range("C15").Value = 'Present the main optical techniques for measuring quantities such as length, displacement and shape, with emphasis on laser interferometric techniques.'
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16 ---
$ws.Rows.Item(16).Clear()
$ws.Range("A3:C3").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = 'Teoria eletromagnética da luz: noções de representação matemática da onda de luz e interpretação de fenômenos como polarização, interferência e difração. Refração, reflexão e óptica geométrica: leis de Snell, equações de Fraunhofer, reflexão total e óptica geométrica. Propagação da luz em meios especiais como cristais fibras ópticas. Óptica de Fourier e holografia: transformada de Fourier e a sua aplicação na óptica como caso de filtros especiais e halográfia. Fontes e sensores de luz: definição e descrição de fontes incoerentes e coerentes e descrição de sensores do tipo puntual, de posição e de imagem. Componentes ópticos e ajuste de sistemas ópticos. Medição de comprimento: método como interferometria, franjas de Moirè, métodos para medição de grandes distâncias. Medição de forma: diversos métodos e técnicas para medição de forma geométrica. Medição de deslocamento, deformação e vibração: métodos de medição que empregam a holografia, speckle" e as franjas de Moirè. Medição de velocidade: métodos de medição de velocidade e sensor de fibras ópticas. Inspeção de falhas: métodos para inspeção de falhas geométricas e internas utilizando a difração ou a difusão da luz.'
$ws.Range("C16").Value = 'Teoria eletromagnética da luz: noções de representação matemática da onda de luz e interpretação de fenômenos como polarização, interferência e difração. Refração, reflexão e óptica geométrica: leis de Snell, equações de Fraunhofer, reflexão total e óptica geométrica. Propagação da luz em meios especiais como cristais fibras ópticas. Óptica de Fourier e holografia: transformada de Fourier e a sua aplicação na óptica como caso de filtros especiais e halográfia. Fontes e sensores de luz: definição e descrição de fontes incoerentes e coerentes e descrição de sensores do tipo puntual, de posição e de imagem. Componentes ópticos e ajuste de sistemas ópticos. Medição de comprimento: método como interferometria, franjas de Moirè, métodos para medição de grandes distâncias. Medição de forma: diversos métodos e técnicas para medição de forma geométrica. Medição de deslocamento, deformação e vibração: métodos de medição que empregam a holografia, speckle" e as franjas de Moirè. Medição de velocidade: métodos de medição de velocidade e sensor de fibras ópticas. Inspeção de falhas: métodos para inspeção de falhas geométricas e internas utilizando a difração ou a difusão da luz.'
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17 ---
$ws.Rows.Item(17).Clear()
$ws.Range("A3:C3").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = 'Electromagnetic theory of light: notions of mathematical representation of the light wave and interpretation of phenomena such as polarization, interference and diffraction. Refraction, reflection and geometric optics: Snell''s laws, Fraunhofer equations, total reflection and geometric optics. Propagation of light in special media such as fiber optic crystals. Fourier optics and holography: Fourier transform and its application in optics as a case of special filters and halography. Light sources and sensors: definition and description of incoherent and coherent sources and description of point, position and image sensors. Optical components and tuning of optical systems. Length measurement: method such as interferometry, Moirè fringes, methods for measuring large distances. Shape measurement: various methods and techniques for geometric shape measurement. Displacement, deformation and vibration measurement: measurement methods employing holography, speckle" and Moirè fringes. Velocity measurement: speed measurement methods and optical fiber sensors. Fault inspection: methods for inspecting geometric and internal spaces using diffraction or scattering of light.'
$ws.Range("C17").Value = 'Electromagnetic theory of light: notions of mathematical representation of the light wave and interpretation of phenomena such as polarization, interference and diffraction. Refraction, reflection and geometric optics: Snell''s laws, Fraunhofer equations, total reflection and geometric optics. Propagation of light in special media such as fiber optic crystals. Fourier optics and holography: Fourier transform and its application in optics as a case of special filters and halography. Light sources and sensors: definition and description of incoherent and coherent sources and description of point, position and image sensors. Optical components and tuning of optical systems. Length measurement: method such as interferometry, Moirè fringes, methods for measuring large distances. Shape measurement: various methods and techniques for geometric shape measurement. Displacement, deformation and vibration measurement: measurement methods employing holography, speckle" and Moirè fringes. Velocity measurement: speed measurement methods and optical fiber sensors. Fault inspection: methods for inspecting geometric and internal spaces using diffraction or scattering of light.'
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18 ---
$ws.Rows.Item(18).Clear()
$ws.Range("A3").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 'Avaliação:'
$ws.Rows.Item(18).AutoFit()

# --- Row 19 ---
$ws.Rows.Item(19).Clear()
$ws.Range("A3:C3").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios.'
$ws.Range("C19").Value = 'Listas de exercícios, provas escritas, apresentação de seminário, aulas de laboratório e preparação de relatórios.'
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20 ---
$ws.Rows.Item(20).Clear()
$ws.Range("A3:C3").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4'
$ws.Range("C20").Value = 'Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4'
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21 ---
$ws.Rows.Item(21).Clear()
$ws.Range("A3:C3").Copy()
$ws.Range("A21:C21").PasteSpecial(-4122)
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 ---
$ws.Rows.Item(22).Clear()
$ws.Range("A3:C3").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'YOSHIZAWA, T. Handbook of Optical Metrology, Boca Raton: CRC Press, 2009.
SALEH, B. E. A.; TEICH, M. C. Handbook of Fotonics, Wiley-Interscience, 2007.
JENKINS, F. A.; WHITE, H. E. Fundamentals of Optics, McGraw-Hill, 1981. 
CREATH, H.; WYANT, J. Measurement of ultraprecision components using non-contact interferometry based instrumentation, Ultraprecision in Manufacturing Engineering, Springer Verlag, 1988.'
$ws.Range("C22").Value = 'YOSHIZAWA, T. Handbook of Optical Metrology, Boca Raton: CRC Press, 2009.
SALEH, B. E. A.; TEICH, M. C. Handbook of Fotonics, Wiley-Interscience, 2007.
JENKINS, F. A.; WHITE, H. E. Fundamentals of Optics, McGraw-Hill, 1981. 
CREATH, H.; WYANT, J. Measurement of ultraprecision components using non-contact interferometry based instrumentation, Ultraprecision in Manufacturing Engineering, Springer Verlag, 1988.'
$ws.Rows.Item(22).RowHeight = 120

# --- Row 23 ---
$ws.Rows.Item(23).Clear()
$ws.Range("A3").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = 'Requisitos:'
$ws.Rows.Item(23).AutoFit()

# --- Row 24 ---
$ws.Rows.Item(24).Clear()
$ws.Range("B3:C3").Copy()
$ws.Range("B24:C24").PasteSpecial(-4122)
$ws.Range("B24").Value = 'LOM3234 -  Óptica Física  (Requisito)
'
$ws.Range("C24").Value = 'LOM3234 -  Óptica Física  (Requisito)
'
$ws.Rows.Item(24).RowHeight = 30

"Edit complete"